$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 70, shifting existing rows 70-132 down to 71-133.
$ws.Rows(70).Insert()

# Populate the new row 70 with a new weekly price observation (same
# Mercado/Producto/Categoria/Variedad/Calidad/Unidad/Origen metadata as the
# rest of the Membrillo "Primera" series), carrying its own date/volume/price
# figures.
$ws.Cells.Item(70, 1).Value = 10
$ws.Cells.Item(70, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(70, 3).Value = "La Araucanía"
$ws.Cells.Item(70, 4).Value = 44638
$ws.Cells.Item(70, 5).Value = 9
$ws.Cells.Item(70, 6).Value = "Fruta"
$ws.Cells.Item(70, 7).Value = 100104
$ws.Cells.Item(70, 8).Value = "Frutos de pepita"
$ws.Cells.Item(70, 9).Value = 100104003
$ws.Cells.Item(70, 10).Value = "Membrillo"
$ws.Cells.Item(70, 11).Value = "Champion"
$ws.Cells.Item(70, 12).Value = "Primera"
$ws.Cells.Item(70, 13).Value = 85
$ws.Cells.Item(70, 14).Value = 22000
$ws.Cells.Item(70, 15).Value = 23000
$ws.Cells.Item(70, 16).Value = 22588
$ws.Cells.Item(70, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(70, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(70, 19).Value = 1255
$ws.Cells.Item(70, 20).Value = 18
